# Applies the "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta" edit:
#  - Adds a new "2508" period block (two new rows) duplicating the existing
#    2507 worker rows, right after the current data rows.
#  - Updates the totals (VALOR MORA, Cant. Periodos) and one salary value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two fresh rows right after the existing data rows (16 & 17),
#    i.e. at rows 18 and 19. Excel automatically re-numbers everything
#    below (merged cells, row 22/23 footer block -> 24/25) and keeps their
#    formatting intact.
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()

# ---------------------------------------------------------------------------
# 2. Row 17 was previously the last data row (bottom-border / "closing"
#    style). Now that two more rows follow it, it becomes a normal middle
#    row, so its bottom border must match row 16's (no special bottom rule).
# ---------------------------------------------------------------------------
for ($c = 2; $c -le 10; $c++) {
    $src = $ws.Cells.Item(16, $c)
    $dst = $ws.Cells.Item(17, $c)
    $dst.Borders.Item(7).LineStyle  = $src.Borders.Item(7).LineStyle   # left
    $dst.Borders.Item(8).LineStyle  = $src.Borders.Item(8).LineStyle   # top
    $dst.Borders.Item(9).LineStyle  = $src.Borders.Item(9).LineStyle   # bottom
    $dst.Borders.Item(10).LineStyle = $src.Borders.Item(10).LineStyle  # right
}

# ---------------------------------------------------------------------------
# 3. Populate the new rows (18 = "middle" style like row 16, 19 = "closing"
#    style like the old row 17) with the new period "2508" duplicating the
#    same two workers / amounts as period "2507".
# ---------------------------------------------------------------------------

# Row 18 - worker ANTONIO (same layout/style as row 16), period 2508
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73158923"
$ws.Range("D18").Value = "ANTONIO JOSE ELJACH GOMEZ"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

for ($c = 2; $c -le 10; $c++) {
    $src = $ws.Cells.Item(16, $c)
    $dst = $ws.Cells.Item(18, $c)
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold
    $dst.NumberFormat = $src.NumberFormat
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment
    $dst.Borders.Item(7).LineStyle  = $src.Borders.Item(7).LineStyle
    $dst.Borders.Item(8).LineStyle  = $src.Borders.Item(8).LineStyle
    $dst.Borders.Item(9).LineStyle  = $src.Borders.Item(9).LineStyle
    $dst.Borders.Item(10).LineStyle = $src.Borders.Item(10).LineStyle
}

# Row 19 - worker ANIBAL (same layout/style as the old closing row 17), period 2508
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73180800"
$ws.Range("D19").Value = "ANIBAL JOSE VERGARA TOUS"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 96000
$ws.Range("G19").Value = 2400000

for ($c = 2; $c -le 10; $c++) {
    $src = $ws.Cells.Item(16, $c)
    $dst = $ws.Cells.Item(19, $c)
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold
    $dst.NumberFormat = $src.NumberFormat
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment
    $dst.Borders.Item(7).LineStyle  = 1
    $dst.Borders.Item(8).LineStyle  = $src.Borders.Item(8).LineStyle
    $dst.Borders.Item(9).LineStyle  = 1
    $dst.Borders.Item(10).LineStyle = 1
}

# ---------------------------------------------------------------------------
# 4. Update the summary header values.
# ---------------------------------------------------------------------------
$ws.Range("G16").Value = 1423500       # Salario Basico, row16 (period 2507 / ANTONIO)
$ws.Range("E11").Value = 305880        # VALOR MORA total (was 152940)
$ws.Range("F13").Value = 2             # Cant. Periodos (was 1)
